# Regenerate orders with updated distance/sizes.
# The Distance codes and one Size code were renumbered:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# These tokens appear (as whole or as part of composite strings like
# "Face04_D64_S20" or "Face04_D64_S20_l.png") in the Condition,
# Filename_Left, Filename_Right, Distance and Size columns. Every other
# column (Trial, Duration_Seconds, Is_Repeat, Block, Face, ConditionID)
# is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $newVal = $val -replace "D64", "D69" -replace "D80", "D86" -replace "D51", "D55" -replace "S30", "S31"
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
